# Applies the Thailand Premier League data update commit.
# - Rows 15/16, 85/86, 117/118: two match records each had swapped rows
#   (match id, teams, score, odds); values below restore the correct pairing.
# - Various other rows: Home/Away team text corrected (Sukhothai FC /
#   Buriram United mixup) -- net text is unchanged where the swap was a pure
#   index fix, so only rows whose actual content differs are written here.
# - Rows 216/217/220/221/223: odds (Q/R/T/U/M/O) corrected for upcoming fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("B15").Value = 6992550
$ws.Range("E15").Value = "Buriram United"
$ws.Range("F15").Value = "Lamphun Warrior FC"
$ws.Range("G15").Value = 3
$ws.Range("I15").Value = "H"
$ws.Range("J15").Value = 1.166
$ws.Range("K15").Value = 8
$ws.Range("L15").Value = 12
$ws.Range("M15").Value = 1.25
$ws.Range("N15").Value = 6.5
$ws.Range("O15").Value = 8.5
$ws.Range("P15").Value = -1.75
$ws.Range("Q15").Value = 1.875
$ws.Range("R15").Value = 1.925
$ws.Range("S15").Value = 3
$ws.Range("V15").Value = 0.25
$ws.Range("W15").Value = -1
$ws.Range("Y15").Value = 0.875
$ws.Range("Z15").Value = -1
$ws.Range("AA15").Value = 0
$ws.Range("AB15").Value = 0

# Row 16
$ws.Range("B16").Value = 6992554
$ws.Range("E16").Value = "Sukhothai FC"
$ws.Range("F16").Value = "Trat FC"
$ws.Range("G16").Value = 0
$ws.Range("I16").Value = "D"
$ws.Range("J16").Value = 1.8
$ws.Range("K16").Value = 3.6
$ws.Range("L16").Value = 4.333
$ws.Range("M16").Value = 1.833
$ws.Range("N16").Value = 3.75
$ws.Range("O16").Value = 4
$ws.Range("P16").Value = -0.5
$ws.Range("Q16").Value = 1.8
$ws.Range("R16").Value = 2
$ws.Range("S16").Value = 2.75
$ws.Range("V16").Value = -1
$ws.Range("W16").Value = 2.75
$ws.Range("Y16").Value = -1
$ws.Range("Z16").Value = 1
$ws.Range("AA16").Value = -1
$ws.Range("AB16").Value = 0.825

# Row 85
$ws.Range("B85").Value = 6992623
$ws.Range("E85").Value = "Ratchaburi FC"
$ws.Range("F85").Value = "Chiangrai Utd"
$ws.Range("G85").Value = 3
$ws.Range("I85").Value = "H"
$ws.Range("J85").Value = 1.7
$ws.Range("K85").Value = 3.75
$ws.Range("L85").Value = 4.2
$ws.Range("M85").Value = 1.7
$ws.Range("N85").Value = 3.75
$ws.Range("O85").Value = 4.333
$ws.Range("P85").Value = -0.75
$ws.Range("Q85").Value = 1.925
$ws.Range("R85").Value = 1.875
$ws.Range("S85").Value = 2.5
$ws.Range("T85").Value = 1.85
$ws.Range("U85").Value = 1.95
$ws.Range("V85").Value = 0.7
$ws.Range("W85").Value = -1
$ws.Range("Y85").Value = 0.925
$ws.Range("Z85").Value = -1
$ws.Range("AA85").Value = 0.8500000000000001
$ws.Range("AB85").Value = -1

# Row 86
$ws.Range("B86").Value = 6992620
$ws.Range("E86").Value = "Uthai Thani FC"
$ws.Range("F86").Value = "Sukhothai FC"
$ws.Range("G86").Value = 0
$ws.Range("I86").Value = "D"
$ws.Range("J86").Value = 1.95
$ws.Range("K86").Value = 3.5
$ws.Range("L86").Value = 3.4
$ws.Range("M86").Value = 2.1
$ws.Range("N86").Value = 3.4
$ws.Range("O86").Value = 3
$ws.Range("P86").Value = -0.25
$ws.Range("Q86").Value = 1.875
$ws.Range("R86").Value = 1.925
$ws.Range("S86").Value = 2.75
$ws.Range("T86").Value = 1.8
$ws.Range("U86").Value = 2
$ws.Range("V86").Value = -1
$ws.Range("W86").Value = 2.4
$ws.Range("Y86").Value = -0.5
$ws.Range("Z86").Value = 0.4625
$ws.Range("AA86").Value = -1
$ws.Range("AB86").Value = 1

# Row 117
$ws.Range("B117").Value = 7485127
$ws.Range("E117").Value = "BG Pathum United"
$ws.Range("F117").Value = "Chiangrai Utd"
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 2
$ws.Range("J117").Value = 1.5
$ws.Range("K117").Value = 4
$ws.Range("L117").Value = 5.75
$ws.Range("M117").Value = 1.363
$ws.Range("N117").Value = 4.5
$ws.Range("O117").Value = 6.5
$ws.Range("P117").Value = -1.25
$ws.Range("Q117").Value = 1.85
$ws.Range("R117").Value = 1.95
$ws.Range("T117").Value = 1.825
$ws.Range("U117").Value = 1.975
$ws.Range("W117").Value = 3.5
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = 0.95
$ws.Range("AA117").Value = 0.825
$ws.Range("AB117").Value = -1

# Row 118
$ws.Range("B118").Value = 7329293
$ws.Range("E118").Value = "Chonburi"
$ws.Range("F118").Value = "Bangkok United"
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 3.6
$ws.Range("K118").Value = 3.5
$ws.Range("L118").Value = 1.85
$ws.Range("M118").Value = 4.5
$ws.Range("N118").Value = 4
$ws.Range("O118").Value = 1.615
$ws.Range("P118").Value = 0.75
$ws.Range("Q118").Value = 1.975
$ws.Range("R118").Value = 1.825
$ws.Range("T118").Value = 1.85
$ws.Range("U118").Value = 1.95
$ws.Range("W118").Value = 3
$ws.Range("Y118").Value = 0.9750000000000001
$ws.Range("Z118").Value = -1
$ws.Range("AA118").Value = -1
$ws.Range("AB118").Value = 0.95

# Row 216
$ws.Range("Q216").Value = 1.825
$ws.Range("R216").Value = 1.975
$ws.Range("T216").Value = 1.85
$ws.Range("U216").Value = 1.95

# Row 217
$ws.Range("M217").Value = 2.05
$ws.Range("O217").Value = 2.875
$ws.Range("Q217").Value = 1.85
$ws.Range("R217").Value = 1.95
$ws.Range("T217").Value = 1.8
$ws.Range("U217").Value = 2

# Row 220
$ws.Range("M220").Value = 2.05
$ws.Range("O220").Value = 3
$ws.Range("Q220").Value = 1.85
$ws.Range("R220").Value = 1.95

# Row 221
$ws.Range("Q221").Value = 1.8
$ws.Range("R221").Value = 2

# Row 223
$ws.Range("Q223").Value = 1.925
$ws.Range("R223").Value = 1.875
